$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'250.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.578"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05680"
$ws.Range("D5").Style = "Normal"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'6.436"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8055"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6MXTokenMX"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.040"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01162"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07288"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03129"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02920"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09273"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001681"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.204"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04722"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006455"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005068"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001051"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.986"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "GateToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D23").Value = "'3.376"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22GateTokenGT"
$ws.Range("D24").Value = "'2.086"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Value = "'0.0003100"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04135"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006880"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1043"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008526"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005631"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.7852"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "'0.01682"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
